$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.327753067016602
$ws.Range("B1").Value = 2.57144570350647
$ws.Range("C1").Value = 2.609463691711426
$ws.Range("D1").Value = 3.312258243560791
$ws.Range("E1").Value = 2.14574146270752
